$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores figures as plain text (e.g. "29.492.09",
# "0.4828") rather than as numbers. Force the cells we are about to
# rewrite to Text format first so Excel keeps our strings exactly as
# typed instead of re-parsing them into floating point numbers.
$priceRows = @(2,3,5,7,8,9,10,11,12,13,14,15,16,17,18,21,22,23,24,25,26,27,28,30,31,32,34,35,36,37,38,39,40,41,44,45,46,47,49,50,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "29.492.09"
$ws.Range("E2").Value = "  +0.69%  "

$ws.Range("D3").Value = "1.913.93"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.68%  "

$ws.Range("D5").Value = "325.94"
$ws.Range("E5").Value = "  +1.44%  "

$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("D7").Value = "0.4828"
$ws.Range("E7").Value = "  +1.93%  "

$ws.Range("D8").Value = "0.4068"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "0.08155"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10").Value = "1.014"
$ws.Range("E10").Value = "  +1.10%  "

$ws.Range("D11").Value = "23.46"
$ws.Range("E11").Value = "  +4.23%  "

$ws.Range("D12").Value = "1.910.39"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").Value = "6.017"
$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("D14").Value = "7.149"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "90.26"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "0.06792"
$ws.Range("E16").Value = "  +2.33%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.69%  "

$ws.Range("D18").Value = "0.00001040"
$ws.Range("E18").Value = "  +1.22%  "

$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("D21").Value = "29.512.67"
$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").Value = "5.629"
$ws.Range("E22").Value = "  +2.07%  "

$ws.Range("D23").Value = "11.76"
$ws.Range("E23").Value = "  +2.67%  "

$ws.Range("D24").Value = "2.184"
$ws.Range("E24").Value = "  -0.66%  "

$ws.Range("D25").Value = "2.140.75"
$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("D26").Value = "155.88"
$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("D27").Value = "6.382"
$ws.Range("E27").Value = "  +6.04%  "

$ws.Range("D28").Value = "20.08"
$ws.Range("E28").Value = "  +1.70%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "119.96"
$ws.Range("E30").Value = "  +2.21%  "

$ws.Range("D31").Value = "1.026"
$ws.Range("E31").Value = "  -4.12%  "

$ws.Range("D32").Value = "0.09529"
$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("D34").Value = "3.562"
$ws.Range("E34").Value = "  +0.34%  "

$ws.Range("D35").Value = "1.390"
$ws.Range("E35").Value = "  -2.35%  "

$ws.Range("D36").Value = "0.02270"
$ws.Range("E36").Value = "  +1.08%  "

$ws.Range("D37").Value = "0.06100"
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").Value = "1.176"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("D39").Value = "0.5968"
$ws.Range("E39").Value = "  +1.93%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "10.81"
$ws.Range("E40").Value = "  +7.09%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "8.012"
$ws.Range("E41").Value = "  -2.77%  "

$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").Value = "2.408"
$ws.Range("E44").Value = "  -4.66%  "

$ws.Range("D45").Value = "12.51"
$ws.Range("E45").Value = "  +3.22%  "

$ws.Range("D46").Value = "0.07617"
$ws.Range("E46").Value = "  -3.48%  "

$ws.Range("D47").Value = "0.5580"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("E48").Value = "  +0.88%  "

$ws.Range("D49").Value = "115.85"
$ws.Range("E49").Value = "  +2.55%  "

$ws.Range("D50").Value = "72.61"
$ws.Range("E50").Value = "  +1.92%  "

$ws.Range("D51").Value = "2.418"
$ws.Range("E51").Value = "  +2.95%  "
